# Worked on type parsing.
# Insert a new "Parameterized?" column into the KgSimpleType grid (column AD),
# shifting the existing "Declaration Scopes" block (old AE:AI) one column to
# the right (new AF:AJ).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column at AD; Excel shifts AD:AI -> AE:AJ and copies the
# adjacent (AC) formatting into the new column for us.
$ws.Columns("AD").Insert()

# Match the width of the rest of the "P:AC" block the new column extends.
$ws.Columns("AD").ColumnWidth = $ws.Columns("AC").ColumnWidth

# Header label for the new column.
$ws.Range("AD1").Value = "Parameterized?"

# Two rows had no existing neighbour in column AC, so Insert() didn't create
# an AD cell for them at all -- give them the same "centered/no-border" style
# used throughout this table (style index 4, same as the rest of column AD).
$ws.Range("AD34").Copy()
$ws.Range("AD33").PasteSpecial(-4122)
$ws.Range("AD44").Copy()
$ws.Range("AD43").PasteSpecial(-4122)

# Fill in the new column's values per row.
$ws.Range("AD10").Value = "Yes"
$ws.Range("AD11").Value = """"
$ws.Range("AD12").Value = """"
$ws.Range("AD13").Value = """"
$ws.Range("AD14").Value = """"
$ws.Range("AD15").Value = """"
$ws.Range("AD16").Value = """"
$ws.Range("AD17").Value = "No"
$ws.Range("AD18").Value = """"
$ws.Range("AD19").Value = """"
$ws.Range("AD20").Value = "No"
$ws.Range("AD21").Value = """"
$ws.Range("AD22").Value = """"
$ws.Range("AD23").Value = """"
$ws.Range("AD24").Value = """"
$ws.Range("AD25").Value = """"
$ws.Range("AD26").Value = """"
$ws.Range("AD27").Value = """"
$ws.Range("AD28").Value = "No"
$ws.Range("AD29").Value = "Yes"
$ws.Range("AD30").Value = "No"
$ws.Range("AD31").Value = "No"
$ws.Range("AD33").Value = "No"
$ws.Range("AD34").Value = """"
$ws.Range("AD35").Value = """"
$ws.Range("AD36").Value = """"
$ws.Range("AD37").Value = """"
$ws.Range("AD38").Value = """"
$ws.Range("AD39").Value = """"
$ws.Range("AD40").Value = """"
$ws.Range("AD41").Value = """"
$ws.Range("AD42").Value = """"
$ws.Range("AD43").Value = """"
$ws.Range("AD44").Value = """"
$ws.Range("AD45").Value = """"
$ws.Range("AD46").Value = """"
$ws.Range("AD47").Value = """"
$ws.Range("AD48").Value = """"
$ws.Range("AD49").Value = """"
$ws.Range("AD50").Value = """"
$ws.Range("AD51").Value = """"
$ws.Range("AD52").Value = """"
$ws.Range("AD53").Value = """"
$ws.Range("AD54").Value = """"
$ws.Range("AD58").Value = "No"
$ws.Range("AD60").Value = "No"
$ws.Range("AD62").Value = "No"
$ws.Range("AD64").Value = "No"
$ws.Range("AD66").Value = "No"
$ws.Range("AD67").Value = """"
$ws.Range("AD68").Value = """"
$ws.Range("AD70").Value = "Yes"
$ws.Range("AD72").Value = "No"

# Insert() leaves behind a few empty, formatted-only cells (AD5, AD7, AD9)
# that weren't part of the edited table -- clear them so no stray cell is
# written for those rows.
$ws.Range("AD5").Clear()
$ws.Range("AD7").Clear()
$ws.Range("AD9").Clear()

# Leave the freshly edited cell selected, matching where the author's cursor
# ended up.
$ws.Range("AD70").Select()
